$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct the PN of the audio cable (row 13, ACCESSORIES section) ---
# Copy the formatting already used by the other "Distributor" column cells
# (e.g. G3) onto G13 so the distributor cell matches the rest of the table.
$ws.Range("G3").Copy()
$ws.Range("G13").PasteSpecial(-4122)

# Distributor PN: was the numeric Farnell part number 3712278, now Digikey's TL1621-ND
$ws.Range("H13").Value = "TL1621-ND"
# MPN: was AV13646, now P312-001
$ws.Range("E13").Value = "P312-001"
# Distributor: was Farnell, now Digikey
$ws.Range("G13").Value = "Digikey"

# --- Update the active cell selection left on the sheet ---
$ws.Range("H23").Select()
